$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'28.314.25"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  -0.66%  "
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(3, 4).Value = "'1.827.59"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  +2.07%  "
$ws.Cells.Item(3, 5).Style = "Normal"

$ws.Cells.Item(4, 4).Value = "'1.003"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "'  -0.27%  "
$ws.Cells.Item(4, 5).Style = "Normal"

$ws.Cells.Item(5, 4).Value = "'329.63"
$ws.Cells.Item(5, 4).Style = "Normal"

$ws.Cells.Item(6, 4).Value = "'1.000"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  -0.06%  "
$ws.Cells.Item(6, 5).Style = "Normal"

$ws.Cells.Item(7, 5).Value = "'  +1.17%  "
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(8, 4).Value = "'0.3761"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "'  +0.09%  "
$ws.Cells.Item(8, 5).Style = "Normal"

$ws.Cells.Item(9, 4).Value = "'44.71"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "'  -2.13%  "
$ws.Cells.Item(9, 5).Style = "Normal"

$ws.Cells.Item(10, 4).Value = "'0.07726"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'  +1.22%  "
$ws.Cells.Item(10, 5).Style = "Normal"

$ws.Cells.Item(11, 4).Value = "'1.130"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "'  -0.91%  "
$ws.Cells.Item(11, 5).Style = "Normal"

$ws.Cells.Item(12, 2).Value = "'BinanceUSD"
$ws.Cells.Item(12, 2).Style = "Normal"
$ws.Cells.Item(12, 3).Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(12, 3).Style = "Normal"
$ws.Cells.Item(12, 4).Value = "'1.000"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "'  -0.47%  "
$ws.Cells.Item(12, 5).Style = "Normal"

$ws.Cells.Item(13, 2).Value = "'Solana"
$ws.Cells.Item(13, 2).Style = "Normal"
$ws.Cells.Item(13, 3).Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(13, 3).Style = "Normal"
$ws.Cells.Item(13, 4).Value = "'22.21"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  -1.94%  "
$ws.Cells.Item(13, 5).Style = "Normal"

$ws.Cells.Item(14, 4).Value = "'6.325"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  +1.00%  "
$ws.Cells.Item(14, 5).Style = "Normal"

$ws.Cells.Item(15, 4).Value = "'7.562"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  +1.35%  "
$ws.Cells.Item(15, 5).Style = "Normal"

$ws.Cells.Item(16, 4).Value = "'1.836.41"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  +2.42%  "
$ws.Cells.Item(16, 5).Style = "Normal"

$ws.Cells.Item(17, 4).Value = "'93.76"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'  +14.69%  "
$ws.Cells.Item(17, 5).Style = "Normal"

$ws.Cells.Item(18, 5).Value = "'  -0.81%  "
$ws.Cells.Item(18, 5).Style = "Normal"

$ws.Cells.Item(19, 4).Value = "'0.06483"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  -3.52%  "
$ws.Cells.Item(19, 5).Style = "Normal"

$ws.Cells.Item(20, 4).Value = "'0.9997"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'  -0.19%  "
$ws.Cells.Item(20, 5).Style = "Normal"

$ws.Cells.Item(21, 4).Value = "'17.53"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  -0.25%  "
$ws.Cells.Item(21, 5).Style = "Normal"

$ws.Cells.Item(22, 4).Value = "'6.332"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  +1.37%  "
$ws.Cells.Item(22, 5).Style = "Normal"

$ws.Cells.Item(23, 4).Value = "'0.5420"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'  +0.61%  "
$ws.Cells.Item(23, 5).Style = "Normal"

$ws.Cells.Item(24, 4).Value = "'28.381.36"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'  -0.61%  "
$ws.Cells.Item(24, 5).Style = "Normal"

$ws.Cells.Item(25, 4).Value = "'11.71"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  -0.98%  "
$ws.Cells.Item(25, 5).Style = "Normal"

$ws.Cells.Item(26, 4).Value = "'2.158"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "'  -11.78%  "
$ws.Cells.Item(26, 5).Style = "Normal"

$ws.Cells.Item(27, 4).Value = "'20.73"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'  +0.59%  "
$ws.Cells.Item(27, 5).Style = "Normal"

$ws.Cells.Item(28, 4).Value = "'155.86"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  +1.54%  "
$ws.Cells.Item(28, 5).Style = "Normal"

$ws.Cells.Item(29, 4).Value = "'2.358"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "'  -0.86%  "
$ws.Cells.Item(29, 5).Style = "Normal"

$ws.Cells.Item(30, 4).Value = "'2.043.91"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "'  +2.38%  "
$ws.Cells.Item(30, 5).Style = "Normal"

$ws.Cells.Item(31, 4).Value = "'128.38"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "'  -1.81%  "
$ws.Cells.Item(31, 5).Style = "Normal"

$ws.Cells.Item(32, 4).Value = "'1.196"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  -8.95%  "
$ws.Cells.Item(32, 5).Style = "Normal"

$ws.Cells.Item(33, 4).Value = "'5.902"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "'  +0.04%  "
$ws.Cells.Item(33, 5).Style = "Normal"

$ws.Cells.Item(34, 4).Value = "'0.09286"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  -0.09%  "
$ws.Cells.Item(34, 5).Style = "Normal"

$ws.Cells.Item(35, 4).Value = "'3.675"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "'  -7.44%  "
$ws.Cells.Item(35, 5).Style = "Normal"

$ws.Cells.Item(36, 4).Value = "'13.15"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "'  +7.81%  "
$ws.Cells.Item(36, 5).Style = "Normal"

$ws.Cells.Item(37, 4).Value = "'0.02349"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "'  +0.43%  "
$ws.Cells.Item(37, 5).Style = "Normal"

$ws.Cells.Item(38, 4).Value = "'0.2189"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  -2.75%  "
$ws.Cells.Item(38, 5).Style = "Normal"

$ws.Cells.Item(39, 4).Value = "'5.186"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  -0.87%  "
$ws.Cells.Item(39, 5).Style = "Normal"

$ws.Cells.Item(40, 4).Value = "'0.6600"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  -1.48%  "
$ws.Cells.Item(40, 5).Style = "Normal"

$ws.Cells.Item(41, 4).Value = "'0.06218"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'  -1.79%  "
$ws.Cells.Item(41, 5).Style = "Normal"

$ws.Cells.Item(42, 4).Value = "'1.205"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'  +0.00%  "
$ws.Cells.Item(42, 5).Style = "Normal"

$ws.Cells.Item(43, 4).Value = "'8.131"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'  +0.76%  "
$ws.Cells.Item(43, 5).Style = "Normal"

$ws.Cells.Item(44, 4).Value = "'0.9994"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  -0.19%  "
$ws.Cells.Item(44, 5).Style = "Normal"

$ws.Cells.Item(45, 5).Value = "'  -1.15%  "
$ws.Cells.Item(45, 5).Style = "Normal"

$ws.Cells.Item(46, 4).Value = "'1.389"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'  -3.94%  "
$ws.Cells.Item(46, 5).Style = "Normal"

$ws.Cells.Item(47, 4).Value = "'0.6103"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  -0.68%  "
$ws.Cells.Item(47, 5).Style = "Normal"

$ws.Cells.Item(48, 4).Value = "'3.777"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  -0.91%  "
$ws.Cells.Item(48, 5).Style = "Normal"

$ws.Cells.Item(49, 4).Value = "'2.051"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'  +1.06%  "
$ws.Cells.Item(49, 5).Style = "Normal"

$ws.Cells.Item(50, 4).Value = "'127.29"
$ws.Cells.Item(50, 4).Style = "Normal"

$ws.Cells.Item(51, 5).Value = "'  +1.28%  "
$ws.Cells.Item(51, 5).Style = "Normal"
